# feat: add 2022-Q1 data
#
# The former "总计" (Total) sheet becomes the new "2022-Q1" per-quarter
# holdings sheet, and a fresh "总计" sheet is appended after it, carrying
# the old totals table plus one new row summarizing 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Turn the existing "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Header row (copy formatting - bold/border/center - from the previous
# quarter's sheet so the new sheet matches the established look).
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$prevQuarter.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Rows = @(
    @(0,  "000654", "华商新锐产业灵活配置混合",         "13.02", "81.79", "3.16",  "0.4114", 4),
    @(1,  "006527", "富国优质发展混合A",                 "13.71", "79.35", "2.69",  "0.3688", 10),
    @(2,  "011212", "富国稳健策略6个月持有期混合A",      "12.61", "83.04", "2.73",  "0.3443", 10),
    @(3,  "004423", "华商研究精选灵活配置混合",         "9.97",  "82.17", "3.17",  "0.3160", 4),
    @(4,  "012491", "华商核心引力混合型证券投资基金A",  "5.05",  "83.13", "3.17",  "0.1601", 4),
    @(5,  "008961", "华商科技创新混合",                 "2.87",  "88.86", "3.91",  "0.1122", 2),
    @(6,  "006528", "富国优质发展混合C",                 "4.04",  "79.35", "2.69",  "0.1087", 10),
    @(7,  "002289", "华商改革创新股票",                 "1.14",  "90.69", "3.69",  "0.0421", 4),
    @(8,  "011213", "富国稳健策略6个月持有期混合C",      "1.14",  "83.04", "2.73",  "0.0311", 10),
    @(9,  "010403", "华商景气优选混合",                 "0.61",  "84.95", "3.73",  "0.0228", 10),
    @(10, "012492", "华商核心引力混合型证券投资基金C",  "0.51",  "83.13", "3.17",  "0.0162", 4),
    @(11, "005161", "华商上游产业股票",                 "0.36",  "89.02", "3.86",  "0.0139", 2),
    @(12, "410006", "华富策略精选混合",                 "0.11",  "77.80", "2.68",  "0.0029", 8)
)

# Columns C..G (fund name / scale / position / ratio / mkt-value) are
# stored as plain text in the source data (e.g. "13.02") rather than
# numbers, so they are entered with a leading apostrophe to force text
# and then "ClearFormats" drops the resulting quote-prefix style.
$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]

    for ($col = 2; $col -le 7; $col++) {
        $cell = $q1.Cells.Item($r, $col)
        $cell.Value = "'" + $row[$col - 1]
        $cell.ClearFormats()
    }

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# Re-apply the "index column" style (bold/border/center) to column A,
# matching every other per-quarter sheet.
$prevQuarter.Range("A2").Copy()
$q1.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($row in $q1Rows) {
    $q1.Cells.Item($row[0] + 2, 1).Value = $row[0]
}

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet with the refreshed totals table.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$prevQuarter.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 13, 1.95),
    @(1, "2021-Q4", 1,  0.03),
    @(2, "2021-Q3", 11, 2.25),
    @(3, "2021-Q2", 8,  2.02),
    @(4, "2021-Q1", 10, 1.58),
    @(5, "2020-Q4", 5,  0.18)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$prevQuarter.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($row in $totalRows) {
    $total.Cells.Item($row[0] + 2, 1).Value = $row[0]
}
